# Auto-generated Excel COM-interop script
# Applies the 'Horarios actualizados Linea 141 - 146' update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value = "Última actualización: 11:11:33"
$ws.Cells.Item(3,1).Value = "Total filas: 117"
$ws.Cells.Item(64,1).Value = "08:27:16"
$ws.Cells.Item(64,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(64,4).Value = 72
$ws.Cells.Item(65,1).Value = "07:50:33"
$ws.Cells.Item(65,3).Value = "15_ABASTO"
$ws.Cells.Item(65,4).Value = 109
$ws.Cells.Item(92,1).Value = "11:11:33"
$ws.Cells.Item(92,2).Value = "11:11"
$ws.Cells.Item(92,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(92,4).Value = 0
$ws.Cells.Item(93,2).Value = "11:14"
$ws.Cells.Item(93,3).Value = "14_ABASTO"
$ws.Cells.Item(93,4).Value = 111
$ws.Cells.Item(94,1).Value = "09:23:23"
$ws.Cells.Item(94,2).Value = "11:15"
$ws.Cells.Item(94,3).Value = "15X38_ABASTO"
$ws.Cells.Item(94,4).Value = 112
$ws.Cells.Item(95,1).Value = "10:37:52"
$ws.Cells.Item(95,2).Value = "11:25"
$ws.Cells.Item(95,3).Value = "16_SANTA ANA"
$ws.Cells.Item(95,4).Value = 48
$ws.Cells.Item(96,1).Value = "10:05:51"
$ws.Cells.Item(96,2).Value = "11:28"
$ws.Cells.Item(96,4).Value = 83
$ws.Cells.Item(97,1).Value = "10:50:41"
$ws.Cells.Item(97,2).Value = "11:29"
$ws.Cells.Item(97,3).Value = "10_OLMOS"
$ws.Cells.Item(97,4).Value = 39
$ws.Cells.Item(98,2).Value = "11:30"
$ws.Cells.Item(98,3).Value = "215C_EL PATO"
$ws.Cells.Item(98,4).Value = 85
$ws.Cells.Item(99,2).Value = "11:31"
$ws.Cells.Item(99,3).Value = "16_SANTA ANA"
$ws.Cells.Item(99,4).Value = 86
$ws.Cells.Item(100,1).Value = "11:11:33"
$ws.Cells.Item(100,2).Value = "11:31"
$ws.Cells.Item(100,3).Value = "215C_EL PATO"
$ws.Cells.Item(100,4).Value = 20
$ws.Cells.Item(101,2).Value = "11:41"
$ws.Cells.Item(101,3).Value = "215B_EL PATO"
$ws.Cells.Item(101,4).Value = 96
$ws.Cells.Item(102,1).Value = "10:05:51"
$ws.Cells.Item(102,2).Value = "11:45"
$ws.Cells.Item(102,3).Value = "15X38_ABASTO"
$ws.Cells.Item(102,4).Value = 100
$ws.Cells.Item(103,1).Value = "11:11:33"
$ws.Cells.Item(103,2).Value = "11:51"
$ws.Cells.Item(103,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(103,4).Value = 40
$ws.Cells.Item(104,1).Value = "10:05:51"
$ws.Cells.Item(104,2).Value = "11:52"
$ws.Cells.Item(104,3).Value = "225_GOMEZ"
$ws.Cells.Item(104,4).Value = 107
$ws.Cells.Item(105,1).Value = "10:37:52"
$ws.Cells.Item(105,2).Value = "11:53"
$ws.Cells.Item(105,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(105,4).Value = 76
$ws.Cells.Item(106,1).Value = "10:50:41"
$ws.Cells.Item(106,2).Value = "11:53"
$ws.Cells.Item(106,3).Value = "225_GOMEZ"
$ws.Cells.Item(106,4).Value = 63
$ws.Cells.Item(107,1).Value = "10:50:41"
$ws.Cells.Item(107,2).Value = "11:54"
$ws.Cells.Item(107,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(107,4).Value = 64
$ws.Cells.Item(108,1).Value = "10:05:51"
$ws.Cells.Item(108,2).Value = "11:58"
$ws.Cells.Item(108,3).Value = "17_ROMERO"
$ws.Cells.Item(108,4).Value = 113
$ws.Cells.Item(109,2).Value = "12:05"
$ws.Cells.Item(109,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(109,4).Value = 88
$ws.Cells.Item(110,2).Value = "12:10"
$ws.Cells.Item(110,3).Value = "15_ABASTO"
$ws.Cells.Item(110,4).Value = 93
$ws.Cells.Item(111,2).Value = "12:10"
$ws.Cells.Item(111,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(111,4).Value = 93
$ws.Cells.Item(112,2).Value = "12:16"
$ws.Cells.Item(112,3).Value = "10_OLMOS"
$ws.Cells.Item(112,4).Value = 99
$ws.Cells.Item(113,1).Value = "11:11:33"
$ws.Cells.Item(113,2).Value = "12:17"
$ws.Cells.Item(113,3).Value = "10_OLMOS"
$ws.Cells.Item(113,4).Value = 66
$ws.Cells.Item(114,1).Value = "10:37:52"
$ws.Cells.Item(114,2).Value = "12:21"
$ws.Cells.Item(114,3).Value = "215C_EL PATO"
$ws.Cells.Item(114,4).Value = 104
$ws.Cells.Item(115,1).Value = "11:11:33"
$ws.Cells.Item(115,2).Value = "12:22"
$ws.Cells.Item(115,3).Value = "215C_EL PATO"
$ws.Cells.Item(115,4).Value = 71
$ws.Cells.Item(115,5).Value = "LP1912"
$ws.Cells.Item(116,1).Value = "10:37:52"
$ws.Cells.Item(116,2).Value = "12:32"
$ws.Cells.Item(116,3).Value = "14_ABASTO"
$ws.Cells.Item(116,4).Value = 115
$ws.Cells.Item(116,5).Value = "LP1912"
$ws.Cells.Item(117,1).Value = "10:37:52"
$ws.Cells.Item(117,2).Value = "12:34"
$ws.Cells.Item(117,3).Value = "15_ABASTO"
$ws.Cells.Item(117,4).Value = 117
$ws.Cells.Item(117,5).Value = "LP1912"
$ws.Cells.Item(118,1).Value = "11:11:33"
$ws.Cells.Item(118,2).Value = "12:35"
$ws.Cells.Item(118,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(118,4).Value = 84
$ws.Cells.Item(118,5).Value = "LP1912"
$ws.Cells.Item(119,1).Value = "10:50:41"
$ws.Cells.Item(119,2).Value = "12:36"
$ws.Cells.Item(119,3).Value = "27_EL RETIRO"
$ws.Cells.Item(119,4).Value = 106
$ws.Cells.Item(119,5).Value = "LP1912"
$ws.Cells.Item(120,1).Value = "10:50:41"
$ws.Cells.Item(120,2).Value = "12:48"
$ws.Cells.Item(120,3).Value = "16_SANTA ANA"
$ws.Cells.Item(120,4).Value = 118
$ws.Cells.Item(120,5).Value = "LP1912"
$ws.Cells.Item(121,1).Value = "11:11:33"
$ws.Cells.Item(121,2).Value = "12:48"
$ws.Cells.Item(121,3).Value = "15X38_ABASTO"
$ws.Cells.Item(121,4).Value = 97
$ws.Cells.Item(121,5).Value = "LP1912"
$ws.Cells.Item(122,1).Value = "11:11:33"
$ws.Cells.Item(122,2).Value = "13:02"
$ws.Cells.Item(122,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(122,4).Value = 111
$ws.Cells.Item(122,5).Value = "LP1912"


# --- Sheet 2: LP1912-215 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,1).Value = "Última actualización: 11:11:33"
$ws.Cells.Item(3,1).Value = "Total filas: 22"
$ws.Cells.Item(24,1).Value = "11:11:33"
$ws.Cells.Item(24,2).Value = "11:31"
$ws.Cells.Item(24,3).Value = "215C_EL PATO"
$ws.Cells.Item(24,4).Value = 20
$ws.Cells.Item(25,1).Value = "10:05:51"
$ws.Cells.Item(25,2).Value = "11:41"
$ws.Cells.Item(25,3).Value = "215B_EL PATO"
$ws.Cells.Item(25,4).Value = 96
$ws.Cells.Item(26,1).Value = "10:37:52"
$ws.Cells.Item(26,2).Value = "12:21"
$ws.Cells.Item(26,3).Value = "215C_EL PATO"
$ws.Cells.Item(26,4).Value = 104
$ws.Cells.Item(26,5).Value = "LP1912"
$ws.Cells.Item(27,1).Value = "11:11:33"
$ws.Cells.Item(27,2).Value = "12:22"
$ws.Cells.Item(27,3).Value = "215C_EL PATO"
$ws.Cells.Item(27,4).Value = 71
$ws.Cells.Item(27,5).Value = "LP1912"


# --- Sheet 3: 6203-6173 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,1).Value = "Última actualización: 11:11:33"

